# sup-figure-4: reorder the "C" and "D" panel-label textboxes so that
# "C" (TextBox 4, id=5) sits just behind "D" (TextBox 5, id=6) in the
# z-order / shape tree instead of in front of it.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shapeC = $s.Shapes.Item("TextBox 4")

# msoSendBackward = 3: move "C" one step back in the z-order, swapping
# places with "D" (which was directly behind it).
$shapeC.ZOrder(3)
